# Updated cryptos list on Sat Jan  6 17:35:43 UTC 2024 with GitHub Actions
# Applies per-cell Price/Volume(1h) updates (and two coin-row reorders) per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.226.82"
$ws.Range("E2").Value = "  +1.37%  "

# Row 3
$ws.Range("D3").Value = "2.245.37"
$ws.Range("E3").Value = "  +1.08%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").Value = "'307.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.27%  "

# Row 6
$ws.Range("D6").Value = "'96.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.01%  "

# Row 7
$ws.Range("E7").Value = "  +1.47%  "

# Row 8
$ws.Range("E8").Value = "  +0.20%  "

# Row 9
$ws.Range("E9").Value = "  -0.45%  "

# Row 10
$ws.Range("D10").Value = "'35.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.90%  "

# Row 11
$ws.Range("D11").Value = "'0.0815"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.14%  "

# Row 12
$ws.Range("D12").Value = "'7.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.85%  "

# Row 13
$ws.Range("E13").Value = "  +0.38%  "

# Row 14
$ws.Range("D14").Value = "2.586.79"
$ws.Range("E14").Value = "  +1.08%  "

# Row 15
$ws.Range("D15").Value = "2.334.36"
$ws.Range("E15").Value = "  +4.82%  "

# Row 16
$ws.Range("D16").Value = "'0.835"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.49%  "

# Row 17
$ws.Range("D17").Value = "'13.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.45%  "

# Row 18
$ws.Range("D18").Value = "44.056.10"
$ws.Range("E18").Value = "  +1.38%  "

# Row 19
$ws.Range("E19").Value = "  +1.81%  "

# Row 20
$ws.Range("E20").Value = "  +2.61%  "

# Row 21
$ws.Range("D21").Value = "'12.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.83%  "

# Row 22
$ws.Range("D22").Value = "'65.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.98%  "

# Row 23
$ws.Range("D23").Value = "'237.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.99%  "

# Row 24
$ws.Range("D24").Value = "'2.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.51%  "

# Row 25
$ws.Range("E25").Value = "  +0.16%  "

# Row 26
$ws.Range("E26").Value = "  +0.09%  "

# Row 27
$ws.Range("D27").Value = "'10.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.09%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.25%  "

# Row 29
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'38.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.35%  "

# Row 30
$ws.Range("D30").Value = "'5.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.84%  "

# Row 31
$ws.Range("D31").Value = "'20.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.08%  "

# Row 32
$ws.Range("D32").Value = "'152.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.25%  "

# Row 33
$ws.Range("D33").Value = "'0.0800"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.48%  "

# Row 34
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.24%  "

# Row 35
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'3.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.15%  "

# Row 36
$ws.Range("E36").Value = "  +3.30%  "

# Row 37
$ws.Range("E37").Value = "  +0.68%  "

# Row 38
$ws.Range("E38").Value = "  -5.06%  "

# Row 39
$ws.Range("D39").Value = "'3.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.32%  "

# Row 40
$ws.Range("D40").Value = "'14.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.67%  "

# Row 41
$ws.Range("D41").Value = "'3.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.13%  "

# Row 42
$ws.Range("D42").Value = "'0.0299"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.75%  "

# Row 43
$ws.Range("E43").Value = "  +0.24%  "

# Row 44
$ws.Range("D44").Value = "1.759.85"
$ws.Range("E44").Value = "  +4.04%  "

# Row 45
$ws.Range("D45").Value = "'83.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.80%  "

# Row 46
$ws.Range("E46").Value = "  +0.16%  "

# Row 47
$ws.Range("D47").Value = "'100.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.41%  "

# Row 48
$ws.Range("D48").Value = "'4.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.45%  "

# Row 49
$ws.Range("D49").Value = "'8.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.87%  "

# Row 50
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "'54.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.95%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.41%  "

